# "Os 100 primeiros de Teste" - fill in the first 100 label rows of the
# "Teste" sheet with the relevance classification (column B), matching the
# header "Irrelevante 0 / Relevante 1" and a couple of special "1-0" cells.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Treinamento")
$ws2 = $wb.Worksheets.Item("Teste")

# Row -> column B value for rows 1..100 of "Teste" (row 1 is the header).
$values = @(
    @(1, "Irrelevante 0 / Relevante 1"),
    @(2, 0),
    @(3, 1),
    @(4, 0),
    @(5, 1),
    @(6, 1),
    @(7, 1),
    @(8, 0),
    @(9, 0),
    @(10, 0),
    @(11, 0),
    @(12, 1),
    @(13, 1),
    @(14, 1),
    @(15, 0),
    @(16, 1),
    @(17, 0),
    @(18, 0),
    @(19, 1),
    @(20, 0),
    @(21, 0),
    @(22, 0),
    @(23, 1),
    @(24, 0),
    @(25, 0),
    @(26, 0),
    @(27, 0),
    @(28, 0),
    @(29, 0),
    @(30, 0),
    @(31, 1),
    @(32, 0),
    @(33, "1-0"),
    @(34, 0),
    @(35, 0),
    @(36, 0),
    @(37, 0),
    @(38, 0),
    @(39, 0),
    @(40, 1),
    @(41, 1),
    @(42, 0),
    @(43, 1),
    @(44, 0),
    @(45, 0),
    @(46, 0),
    @(47, 0),
    @(48, 0),
    @(49, 0),
    @(50, 0),
    @(51, 0),
    @(52, 0),
    @(53, 1),
    @(54, 0),
    @(55, 1),
    @(56, 0),
    @(57, 0),
    @(58, 0),
    @(59, 0),
    @(60, 1),
    @(61, 1),
    @(62, 0),
    @(63, 0),
    @(64, 0),
    @(65, 0),
    @(66, 0),
    @(67, 0),
    @(68, 0),
    @(69, 1),
    @(70, "1-0"),
    @(71, 1),
    @(72, 1),
    @(73, 0),
    @(74, 0),
    @(75, 1),
    @(76, 0),
    @(77, 0),
    @(78, 0),
    @(79, 0),
    @(80, 0),
    @(81, 0),
    @(82, 0),
    @(83, 1),
    @(84, 0),
    @(85, 0),
    @(86, 1),
    @(87, 0),
    @(88, 0),
    @(89, 0),
    @(90, 1),
    @(91, 0),
    @(92, 0),
    @(93, 1),
    @(94, 0),
    @(95, 0),
    @(96, 0),
    @(97, 1),
    @(98, 0),
    @(99, 1),
    @(100, 1)
)

foreach ($pair in $values) {
    $row = $pair[0]
    $val = $pair[1]
    $ws2.Range("B$row").Value = $val
}

# Column widths on "Teste" (A holds the long tweet text, B the label).
$ws2.Columns.Item(1).ColumnWidth = 139.26
$ws2.Columns.Item(2).ColumnWidth = 23.01

# View state: "Treinamento" loses the active-tab/scroll position it had,
# moving its selection back up to B1; "Teste" becomes the active sheet,
# scrolled near the bottom of the just-entered data with B101 selected
# (the next empty cell to fill) at 75% zoom.
$ws1.Range("B1").Select()

$ws2.Activate()
$ws2.Range("B101").Select()
$excel.ActiveWindow.ScrollRow = 159
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 75
